$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação: date string update (force text so "01/01/2022" is not
# auto-converted into an Excel date serial, matching the original
# plain-text shared-string storage).
$ws.Range("B8:C8").NumberFormat = "@"
$ws.Range("B8:C8").Value = "01/01/2022"

# Programa resumido: (item "3) Determinação de potência de bomba" removed,
# list renumbered, and spacing normalized after each "N)" marker)
$ws.Range("B14:C14").Value = "1) Determinação do número de Reynolds; 2) Medidas de vazão em líquidos; 3) Esvaziamento de tanques cilíndricos; 4) Determinação do fator de atrito em tubulações; 5) Ensaios reológicos."

# Short syllabus: (spacing normalized after "1)")
$ws.Range("B15:C15").Value = "1) Determination of Reynolds number 2) Liquids flow measurements 3) Determination of the fluids velocity profile in pipe 4) Liquid flow from a cylindrical tank 5) Determination of friction factor in pipes 6) Rheological tests."

# Programa: (item 3 about the centrifugal pump removed, list renumbered)
$ws.Range("B16:C16").Value = "1) Determinação do número de Reynolds, utilizando o aparato experimental; 2) Medidas de Vazão em líquidos: determinação da vazão com utilização de placa de orifício, Venturi e rotâmetro; 3) Esvaziamento de tanques cilíndricos: aplicação das equações de conservação de massa e energia na determinação do tempo de esvaziamento de tanques e comparação com dados experimentais; 4) Determinação do fator de atrito em tubulações: avaliação do fator de atrito em função do número de Reynolds em tubulações; 5) Ensaios reológicos: com a utilização de vários tipos de viscosímetros para determinar as viscosidades dinâmica e cinemática."

# Syllabus: (item 3 about the velocity profile / Pitot tube removed, list renumbered)
$ws.Range("B17:C17").Value = "1) Determination of the Reynolds number using the experimental apparatus. 2) Liquids flow measurements: Determining the flow by using a rotameter, orifice plate and Venturi. 3) Liquid flow from a cylindrical tank: verification of a mathematical model based on the conservation equations of mass and energy to determine the time of emptying reservoir and compared the results with experimental data. 4) Determination of friction factor in pipes: evaluation of the friction factor as a function of Reynolds number in pipes. 5) Rheological tests: with the use of various types of viscometers to determine the dynamic and kinematic viscosities."
